# Update the "Förändrad" (Changed) date in column C for rows 2-18
# from serial date 45207 (2023-10-08) to 45208 (2023-10-09).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 18; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45207) {
        $cell.Value2 = 45208
    }
}
